# Generate Report for Handback
# Two source files (b6e907db..., cef68b5a...) have now been handed back
# ("Handed back: in sync with en-US") for both the zh-cn and de-de locales.
# Update the Overview sheet status columns, and on each locale sheet fill
# in the Status / Latest Target File / Latest Handback File / Latest
# Handback DateTime columns for those two rows.

$wb = $excel.ActiveWorkbook

$HANDED_BACK = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (b6e907db...) and 5 (cef68b5a...) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# zh-cn (E) and de-de (F) columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $HANDED_BACK
$wsOverview.Range("F4").Value = $HANDED_BACK
$wsOverview.Range("E5").Value = $HANDED_BACK
$wsOverview.Range("F5").Value = $HANDED_BACK

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 (b6e907db...) and 5 (cef68b5a...).
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $HANDED_BACK
$wsZhCn.Range("I4").Value = "b6e907db-2e4d-40f9-92e2-db72b3e1f51f.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/83b9c2af67e2e4cb15b44d6f68535f43fee7e18f/e2e/b6e907db-2e4d-40f9-92e2-db72b3e1f51f.md", [Type]::Missing, [Type]::Missing, "b6e907db-2e4d-40f9-92e2-db72b3e1f51f.md")
$wsZhCn.Range("I4").Font.Underline = 1
$wsZhCn.Range("I4").Font.Color = 15570276
$wsZhCn.Range("J4").Value = "b6e907db-2e4d-40f9-92e2-db72b3e1f51f.0a3d65b98b33466e8bcc15f9e07d354879b62d0f.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-25 10:31:59"

$wsZhCn.Range("C5").Value = $HANDED_BACK
$wsZhCn.Range("I5").Value = "cef68b5a-7ace-4139-a89e-2771efe58003.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/83b9c2af67e2e4cb15b44d6f68535f43fee7e18f/e2e/cef68b5a-7ace-4139-a89e-2771efe58003.md", [Type]::Missing, [Type]::Missing, "cef68b5a-7ace-4139-a89e-2771efe58003.md")
$wsZhCn.Range("I5").Font.Underline = 1
$wsZhCn.Range("I5").Font.Color = 15570276
$wsZhCn.Range("J5").Value = "cef68b5a-7ace-4139-a89e-2771efe58003.c3ebdc56c6e733995306629c035e04289619ba29.zh-cn.xlf"
$wsZhCn.Range("K5").Value = "2016-08-25 10:31:59"

# ---------------------------------------------------------------------
# de-de sheet: rows 4 (b6e907db...) and 5 (cef68b5a...).
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $HANDED_BACK
$wsDeDe.Range("I4").Value = "b6e907db-2e4d-40f9-92e2-db72b3e1f51f.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c0e4959ffe2678b326e3329e6c30c4d9904cc9aa/e2e/b6e907db-2e4d-40f9-92e2-db72b3e1f51f.md", [Type]::Missing, [Type]::Missing, "b6e907db-2e4d-40f9-92e2-db72b3e1f51f.md")
$wsDeDe.Range("I4").Font.Underline = 1
$wsDeDe.Range("I4").Font.Color = 15570276
$wsDeDe.Range("J4").Value = "b6e907db-2e4d-40f9-92e2-db72b3e1f51f.0a3d65b98b33466e8bcc15f9e07d354879b62d0f.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-25 10:32:15"

$wsDeDe.Range("C5").Value = $HANDED_BACK
$wsDeDe.Range("I5").Value = "cef68b5a-7ace-4139-a89e-2771efe58003.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c0e4959ffe2678b326e3329e6c30c4d9904cc9aa/e2e/cef68b5a-7ace-4139-a89e-2771efe58003.md", [Type]::Missing, [Type]::Missing, "cef68b5a-7ace-4139-a89e-2771efe58003.md")
$wsDeDe.Range("I5").Font.Underline = 1
$wsDeDe.Range("I5").Font.Color = 15570276
$wsDeDe.Range("J5").Value = "cef68b5a-7ace-4139-a89e-2771efe58003.c3ebdc56c6e733995306629c035e04289619ba29.de-de.xlf"
$wsDeDe.Range("K5").Value = "2016-08-25 10:32:15"

Write-Output "localization-status report updated for handback"
